$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the four "additional papers" hyperlink paragraphs that
#    follow the "Literature Review - Additional Papers" heading.
# ------------------------------------------------------------------
$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "nature\.com/articles/s41577-022-00687-3") {
        $startPara = $i
    }
    if ($t -match "mdpi\.com/2076-393X/10/4/591") {
        $endPara = $i
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $rangeStart = $d.Paragraphs.Item($startPara).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endPara).Range.End
    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}

# ------------------------------------------------------------------
# 2) Add the (latent, built-in) "FollowedHyperlink" character style
#    to the style sheet, matching what Word mints once a followed
#    hyperlink style is first referenced.
# ------------------------------------------------------------------
$style = $d.Styles.Add("FollowedHyperlink", 2)
$style.BaseStyle = "DefaultParagraphFont"
$style.Priority = 99
$style.UnhideWhenUsed = $true
$style.Font.Color = 7491477
$style.Font.Underline = 1
